$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.450.52"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").Value = "2.646.99"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.20"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.86"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").Value = "2.645.90"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("E10").Value = "  +7.10%  "
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("E13").Value = "  +2.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.15"
$ws.Range("E14").Value = "  +1.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000192"
$ws.Range("E15").Value = "  +2.39%  "
$ws.Range("D16").Value = "3.128.31"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "68.345.17"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "2.650.39"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.39"
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "364.03"
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.39"
$ws.Range("E22").Value = "  +3.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.90"
$ws.Range("E23").Value = "  +1.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.69"
$ws.Range("E25").Value = "  +3.41%  "
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("E28").Value = "  +1.86%  "
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "574.73"
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("E32").Value = "  +3.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.43"
$ws.Range("E33").Value = "  +1.53%  "
$ws.Range("E34").Value = "  +2.04%  "
$ws.Range("E35").Value = "  +3.29%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  +5.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.02"
$ws.Range("E38").Value = "  +2.02%  "
$ws.Range("E39").Value = "  +0.85%  "
$ws.Range("E40").Value = "  +1.46%  "
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").Value = "0.0₆0337"
$ws.Range("E43").Value = "  -1.81%  "
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.74"
$ws.Range("E45").Value = "  +3.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.68"
$ws.Range("E46").Value = "  +1.08%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "157.16"
$ws.Range("E48").Value = "  +0.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.77"
$ws.Range("E49").Value = "  +2.05%  "
$ws.Range("E50").Value = "  +0.74%  "
